# Refreshes the cryptocurrency price/volume snapshot (GitHub Actions data pull).
# Rows 2-51 map 1:1 onto the sheet rows; B/C only change for the LidoDAOToken /
# TheGraph row swap at rows 42-43.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.774.70"
$ws.Range("E2").Value = "  +1.09%  "

# Row 3
$ws.Range("D3").Value = "3.467.26"
$ws.Range("E3").Value = "  +1.10%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").Value = "413.71"
$ws.Range("E5").Value = "  +0.97%  "

# Row 6
$ws.Range("D6").Value = "130.17"
$ws.Range("E6").Value = "  -0.22%  "

# Row 7
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  -0.94%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").Value = "0.726"
$ws.Range("E9").Value = "  -1.85%  "

# Row 10
$ws.Range("D10").Value = "0.151"
$ws.Range("E10").Value = "  +5.85%  "

# Row 11
$ws.Range("D11").Value = "42.55"
$ws.Range("E11").Value = "  -1.06%  "

# Row 12
$ws.Range("D12").Value = "'9.60"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.00%  "

# Row 13
$ws.Range("E13").Value = "  -2.56%  "

# Row 14
$ws.Range("D14").Value = "4.022.57"
$ws.Range("E14").Value = "  +1.15%  "

# Row 15
$ws.Range("D15").Value = "'0.140"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.29%  "

# Row 16
$ws.Range("E16").Value = "  -3.89%  "

# Row 17
$ws.Range("D17").Value = "3.456.12"
$ws.Range("E17").Value = "  +0.96%  "

# Row 18
$ws.Range("D18").Value = "12.68"
$ws.Range("E18").Value = "  +1.17%  "

# Row 19
$ws.Range("E19").Value = "  -1.66%  "

# Row 20
$ws.Range("D20").Value = "62.681.62"
$ws.Range("E20").Value = "  +1.03%  "

# Row 21
$ws.Range("D21").Value = "461.49"
$ws.Range("E21").Value = "  +1.25%  "

# Row 22
$ws.Range("D22").Value = "90.55"
$ws.Range("E22").Value = "  -1.29%  "

# Row 23
$ws.Range("E23").Value = "  +1.76%  "

# Row 24
$ws.Range("D24").Value = "13.31"
$ws.Range("E24").Value = "  +1.03%  "

# Row 25
$ws.Range("D25").Value = "10.73"
$ws.Range("E25").Value = "  +17.58%  "

# Row 26
$ws.Range("E26").Value = "  +0.51%  "

# Row 27
$ws.Range("D27").Value = "33.42"
$ws.Range("E27").Value = "  +0.54%  "

# Row 28
$ws.Range("D28").Value = "'4.80"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.41%  "

# Row 29
$ws.Range("D29").Value = "7.59"
$ws.Range("E29").Value = "  -1.27%  "

# Row 30
$ws.Range("E30").Value = "  -0.99%  "

# Row 31
$ws.Range("E31").Value = "  -1.26%  "

# Row 32
$ws.Range("E32").Value = "  -2.31%  "

# Row 33
$ws.Range("E33").Value = "  -1.98%  "

# Row 34
$ws.Range("D34").Value = "40.77"
$ws.Range("E34").Value = "  -5.67%  "

# Row 35
$ws.Range("E35").Value = "  +0.02%  "

# Row 36
$ws.Range("D36").Value = "58.47"
$ws.Range("E36").Value = "  +7.64%  "

# Row 37
$ws.Range("D37").Value = "'0.0490"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.61%  "

# Row 38
$ws.Range("E38").Value = "  +4.59%  "

# Row 39
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.16%  "

# Row 40
$ws.Range("D40").Value = "147.47"
$ws.Range("E40").Value = "  +3.57%  "

# Row 41
$ws.Range("E41").Value = "  -0.82%  "

# Row 42
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value = "0.321"
$ws.Range("E42").Value = "  +0.44%  "

# Row 43
$ws.Range("B43").Value = "LidoDAOToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D43").Value = "3.33"
$ws.Range("E43").Value = "  -1.16%  "

# Row 44
$ws.Range("E44").Value = "  +5.28%  "

# Row 45
$ws.Range("E45").Value = "  +2.01%  "

# Row 46
$ws.Range("D46").Value = "2.07"
$ws.Range("E46").Value = "  +3.41%  "

# Row 47
$ws.Range("D47").Value = "2.41"
$ws.Range("E47").Value = "  +13.67%  "

# Row 48
$ws.Range("D48").Value = "0.0₃0557"
$ws.Range("E48").Value = "  +29.87%  "

# Row 49
$ws.Range("E49").Value = "  -1.53%  "

# Row 50
$ws.Range("D50").Value = "22.43"
$ws.Range("E50").Value = "  -0.22%  "

# Row 51
$ws.Range("E51").Value = "  -0.86%  "
